$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 473.86667
$ws.Range("I6").Value = 473.86667
$ws.Range("K6").Value = 1421.60001
$ws.Range("M6").Value = -1309.60001
$ws.Range("H129").Value = 1273.975
$ws.Range("J129").Value = 1376.9166
$ws.Range("L129").Value = 4130.7498
$ws.Range("N129").Value = -14130.7498
$ws.Range("H138").Value = 2357.1177
$ws.Range("I138").Value = 1382.2
$ws.Range("J138").Value = 3749.8572
$ws.Range("K138").Value = 4146.6
$ws.Range("L138").Value = 11249.5716
$ws.Range("M138").Value = 993.3999999999996
$ws.Range("N138").Value = -21529.5716
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2052.5715
$ws.Range("I61").Value = 2092.7058
$ws.Range("J61").Value = 1882
$ws.Range("K61").Value = 2092.7058
$ws.Range("L61").Value = 1882
$ws.Range("M61").Value = -1880.7058
$ws.Range("N61").Value = -2306
$ws.Range("H122").Value = 4268.8335
$ws.Range("J122").Value = 5650.25
$ws.Range("L122").Value = 16950.75
$ws.Range("N122").Value = -21850.75
$ws.Range("H132").Value = 2720.4048
$ws.Range("I132").Value = 2068.4644
$ws.Range("J132").Value = 4024.2856
$ws.Range("K132").Value = 6205.3932
$ws.Range("L132").Value = 12072.8568
$ws.Range("M132").Value = -3675.3932
$ws.Range("N132").Value = -17132.8568
$ws.Range("H136").Value = 2052.5715
$ws.Range("I136").Value = 2092.7058
$ws.Range("J136").Value = 1882
$ws.Range("K136").Value = 6278.117400000001
$ws.Range("L136").Value = 5646
$ws.Range("M136").Value = -3728.117400000001
$ws.Range("N136").Value = -10746
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 206882.88
$ws.Range("I31").Value = 356381.16
$ws.Range("J31").Value = 3992.3572
$ws.Range("K31").Value = 356381.16
$ws.Range("L31").Value = 3992.3572
$ws.Range("M31").Value = -356086.16
$ws.Range("N31").Value = -4582.3572
$ws.Range("H34").Value = 206882.88
$ws.Range("I34").Value = 356381.16
$ws.Range("J34").Value = 3992.3572
$ws.Range("K34").Value = 356381.16
$ws.Range("L34").Value = 3992.3572
$ws.Range("M34").Value = -356179.16
$ws.Range("N34").Value = -4396.3572
$ws.Range("H132").Value = 3173.465
$ws.Range("I132").Value = 2619.6072
$ws.Range("K132").Value = 7858.821599999999
$ws.Range("M132").Value = -5328.821599999999
$ws.Range("H137").Value = 44551.43
$ws.Range("J137").Value = 44551.43
$ws.Range("L137").Value = 44551.43
$ws.Range("N137").Value = -54751.43
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1528.5714
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 1528.5714
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H113").Value = 1761157
$ws.Range("I113").Value = 599.58185
$ws.Range("J113").Value = 7813073
$ws.Range("K113").Value = 1798.74555
$ws.Range("L113").Value = 23439219
$ws.Range("M113").Value = 371.2544499999999
$ws.Range("N113").Value = -23443559
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 14009
$ws.Range("J6").Value = 14009
$ws.Range("L6").Value = 14009
$ws.Range("N6").Value = -14235
$ws.Range("H16").Value = 14009
$ws.Range("J16").Value = 14009
$ws.Range("L16").Value = 14009
$ws.Range("N16").Value = -14509
$ws.Range("H70").Value = 6144.304
$ws.Range("I70").Value = 5503.4863
$ws.Range("J70").Value = 8778.777
$ws.Range("K70").Value = 5503.4863
$ws.Range("L70").Value = 8778.777
$ws.Range("M70").Value = -5233.4863
$ws.Range("N70").Value = -9318.777
$ws.Range("H73").Value = 6144.304
$ws.Range("I73").Value = 5503.4863
$ws.Range("J73").Value = 8778.777
$ws.Range("K73").Value = 5503.4863
$ws.Range("L73").Value = 8778.777
$ws.Range("M73").Value = -4567.4863
$ws.Range("N73").Value = -10650.777
$ws.Range("H102").Value = 4154.6665
$ws.Range("I102").Value = 3232.1
$ws.Range("K102").Value = 3232.1
$ws.Range("M102").Value = -1610.1
$ws.Range("H126").Value = 3747.5144
$ws.Range("I126").Value = 2931
$ws.Range("K126").Value = 8793
$ws.Range("M126").Value = -6323
$ws.Range("H132").Value = 2090.403
$ws.Range("I132").Value = 1716.34
$ws.Range("J132").Value = 3190.5881
$ws.Range("K132").Value = 5149.02
$ws.Range("L132").Value = 9571.764299999999
$ws.Range("M132").Value = -2619.02
$ws.Range("N132").Value = -14631.7643
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4549
$ws.Range("I7").Value = 2023.5
$ws.Range("K7").Value = 2023.5
$ws.Range("M7").Value = -1911.5
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("N20").ClearContents()
$ws.Range("H74").Value = 44062
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 44062
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 44062
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -46058
$ws.Range("H77").Value = 44062
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 44062
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 132186
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -142170
$ws.Range("H126").Value = 4549
$ws.Range("I126").Value = 2023.5
$ws.Range("K126").Value = 6070.5
$ws.Range("M126").Value = -3600.5
$ws.Range("H132").Value = 3378.0625
$ws.Range("I132").Value = 2504.0908
$ws.Range("J132").Value = 5300.8
$ws.Range("K132").Value = 7512.2724
$ws.Range("L132").Value = 15902.4
$ws.Range("M132").Value = -4982.2724
$ws.Range("N132").Value = -20962.4
$ws.Range("H136").Value = 2507.9387
$ws.Range("I136").Value = 971.6389
$ws.Range("J136").Value = 6762.3076
$ws.Range("K136").Value = 2914.9167
$ws.Range("L136").Value = 20286.9228
$ws.Range("M136").Value = -364.9167000000002
$ws.Range("N136").Value = -25386.9228
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 34647.5
$ws.Range("J123").Value = 34647.5
$ws.Range("L123").Value = 34647.5
$ws.Range("N123").Value = -44447.5
$ws.Range("H133").Value = 34839.375
$ws.Range("J133").Value = 34839.375
$ws.Range("L133").Value = 34839.375
$ws.Range("N133").Value = -44959.375
$ws.Range("H136").Value = 2228.3333
$ws.Range("I136").Value = 880.7879
$ws.Range("J136").Value = 4698.8335
$ws.Range("K136").Value = 2642.3637
$ws.Range("L136").Value = 14096.5005
$ws.Range("M136").Value = -92.36369999999988
$ws.Range("N136").Value = -19196.5005
